# Auto-generated edit script applying the Cuchulainn_Profits.xlsx diff
# Updates currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 91.333336
$ws.Range("I5").Value = 91.333336
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 91.333336
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 23.666664
$ws.Range("N5").ClearContents()

# Row 38 (Leve Item ID 4599)
$ws.Range("H38").Value = 2179.9167
$ws.Range("I38").Value = 11.75
$ws.Range("K38").Value = 35.25
$ws.Range("M38").Value = 336.75

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 9077
$ws.Range("I116").Value = 900
$ws.Range("J116").Value = 9985.556
$ws.Range("K116").Value = 900
$ws.Range("L116").Value = 9985.556
$ws.Range("M116").Value = 2542
$ws.Range("N116").Value = -16869.556

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 1413.1538
$ws.Range("I32").Value = 1413.1538
$ws.Range("K32").Value = 1413.1538
$ws.Range("M32").Value = -1126.1538

# Row 62 (Leve Item ID 10719)
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65 (Leve Item ID 10719)
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 29 (Leve Item ID 2318)
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# Row 97 (Leve Item ID 19518)
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

# Row 112 (Leve Item ID 25788)
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# Row 124 (Leve Item ID 34245)
$ws.Range("H124").Value = 15780
$ws.Range("J124").Value = 15780
$ws.Range("L124").Value = 15780
$ws.Range("N124").Value = -25600

# Row 130 (Leve Item ID 34682)
$ws.Range("H130").Value = 42500
$ws.Range("J130").Value = 42500
$ws.Range("L130").Value = 42500
$ws.Range("N130").Value = -52540

# Row 135 (Leve Item ID 41992)
$ws.Range("H135").Value = 59997.5
$ws.Range("J135").Value = 59997.5
$ws.Range("L135").Value = 59997.5
$ws.Range("N135").Value = -70137.5

$ws = $wb.Worksheets.Item("CRP")
# Row 42 (Leve Item ID 1847)
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 6378.5
$ws.Range("I58").Value = 1683.7142
$ws.Range("K58").Value = 1683.7142
$ws.Range("M58").Value = -1480.7142

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 7759.778
$ws.Range("I134").Value = 1535.75
$ws.Range("J134").Value = 12739
$ws.Range("K134").Value = 4607.25
$ws.Range("L134").Value = 38217
$ws.Range("M134").Value = -2072.25
$ws.Range("N134").Value = -43287

# Row 135 (Leve Item ID 42008)
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 6378.5
$ws.Range("I136").Value = 1683.7142
$ws.Range("K136").Value = 5051.142599999999
$ws.Range("M136").Value = -2501.142599999999

# Row 137 (Leve Item ID 43231)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 5188.8696
$ws.Range("I4").Value = 331.81818
$ws.Range("J4").Value = 9641.166999999999
$ws.Range("K4").Value = 995.45454
$ws.Range("L4").Value = 28923.501
$ws.Range("M4").Value = -883.45454
$ws.Range("N4").Value = -29147.501

# Row 48 (Leve Item ID 4724)
$ws.Range("H48").Value = 250
$ws.Range("I48").Value = 250
$ws.Range("K48").Value = 750
$ws.Range("M48").Value = -500

# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()

# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()

# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 1000
$ws.Range("J107").Value = 1000
$ws.Range("L107").Value = 3000
$ws.Range("N107").Value = -6840

# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 349.66666
$ws.Range("I122").Value = 240
$ws.Range("J122").Value = 898
$ws.Range("K122").Value = 2160
$ws.Range("L122").Value = 8082
$ws.Range("M122").Value = 290
$ws.Range("N122").Value = -12982

# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 2876.6
$ws.Range("J129").Value = 2847.25
$ws.Range("L129").Value = 8541.75
$ws.Range("N129").Value = -18541.75

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 2117.9333
$ws.Range("J131").Value = 2576.818
$ws.Range("L131").Value = 7730.454000000001
$ws.Range("N131").Value = -17810.454

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1050
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 8100
$ws.Range("L132").Value = 10800
$ws.Range("M132").Value = -5570
$ws.Range("N132").Value = -15860

# Row 134 (Leve Item ID 44074)
$ws.Range("H134").Value = 5078
$ws.Range("I134").Value = 4695.8
$ws.Range("K134").Value = 14087.4
$ws.Range("M134").Value = -9017.400000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 103 (Leve Item ID 26023)
$ws.Range("H103").Value = 249999
$ws.Range("J103").Value = 249999
$ws.Range("L103").Value = 249999
$ws.Range("N103").Value = -252343

# Row 123 (Leve Item ID 34150)
$ws.Range("H123").Value = 47000
$ws.Range("J123").Value = 47000
$ws.Range("L123").Value = 47000
$ws.Range("N123").Value = -51900

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2499.5
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 8997
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -13937

# Row 141 (Leve Item ID 42504)
$ws.Range("H141").Value = 88570.164
$ws.Range("J141").Value = 88570.164
$ws.Range("L141").Value = 88570.164
$ws.Range("N141").Value = -98930.164

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 140 (Leve Item ID 42506)
$ws.Range("H140").Value = 62997.5
$ws.Range("J140").Value = 62997.5
$ws.Range("L140").Value = 62997.5
$ws.Range("N140").Value = -73357.5
